$d = $word.ActiveDocument

# Update the final exam date range to just the end date
$d.Content.Find.Execute("Mon 12/8 – Thurs 12/11", $true, $false, $false, $false, $false, $true, 1, $false, "Thurs 12/11", 2)

# Update the final exam time from TBD to the actual scheduled time
$d.Content.Find.Execute("TBD", $true, $false, $false, $false, $false, $true, 1, $false, "10:00 – 11:50 am", 2)
